$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - text values, forced via NumberFormat "@" so the
# engine does not auto-coerce numeric-looking strings (e.g. "2.60") into numbers.
$priceUpdates = @{
    2 = "50.072.08"
    3 = "2.661.25"
    5 = "114.29"
    6 = "326.49"
    7 = "0.529"
    8 = "0.999"
    9 = "0.557"
    10 = "41.29"
    11 = "20.16"
    15 = "3.077.02"
    16 = "2.645.06"
    17 = "0.879"
    18 = "50.011.74"
    19 = "13.33"
    20 = "6.81"
    21 = "2.93"
    22 = "0.0₃0960"
    23 = "72.62"
    24 = "278.51"
    25 = "2.60"
    26 = "26.94"
    28 = "10.04"
    30 = "36.77"
    32 = "50.35"
    33 = "5.52"
    35 = "0.0819"
    37 = "5.06"
    39 = "3.13"
    40 = "125.35"
    42 = "22.28"
    45 = "2.109.00"
    49 = "9.12"
    51 = "59.82"
}
foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Column E (Volume 1h) updates - plain text percentages, no coercion risk.
$volumeUpdates = @{
    2 = "  +4.37%  "
    3 = "  +7.16%  "
    4 = "  +0.07%  "
    5 = "  +8.53%  "
    6 = "  +2.92%  "
    7 = "  +2.19%  "
    9 = "  +3.86%  "
    10 = "  +6.36%  "
    11 = "  +0.13%  "
    12 = "  +3.23%  "
    13 = "  +0.28%  "
    14 = "  +4.36%  "
    15 = "  +7.12%  "
    16 = "  +6.30%  "
    17 = "  +6.51%  "
    18 = "  +4.41%  "
    19 = "  +4.92%  "
    20 = "  +3.90%  "
    21 = "  -0.88%  "
    22 = "  +3.54%  "
    23 = "  +2.45%  "
    24 = "  +2.33%  "
    25 = "  +3.86%  "
    26 = "  +4.99%  "
    28 = "  +3.66%  "
    29 = "  +1.71%  "
    30 = "  +6.60%  "
    31 = "  +3.30%  "
    32 = "  +2.00%  "
    33 = "  +5.11%  "
    34 = "  +4.19%  "
    35 = "  +6.41%  "
    36 = "  -0.15%  "
    37 = "  +11.16%  "
    38 = "  +8.12%  "
    39 = "  +9.54%  "
    40 = "  +2.44%  "
    41 = "  +2.47%  "
    42 = "  +0.99%  "
    43 = "  +0.50%  "
    44 = "  +5.68%  "
    45 = "  +5.47%  "
    46 = "  +6.19%  "
    47 = "  +13.55%  "
    48 = "  +5.38%  "
    49 = "  +2.58%  "
    50 = "  +4.16%  "
    51 = "  +6.29%  "
}
foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
